# Generate Report for Handback
# Fills in the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" / "Error Detail" columns for the last row
# (the 65a8a149-... entry) on both the zh-cn and de-de sheets, because a
# handback was produced for that entry but the generated handback file's
# version does not match the latest source, so an error is also recorded.
# Also widens the "Error Detail" column so the long message is readable.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7b20ac077373de7e28dde12370a55d1901557105/e2e/65a8a149-dbf1-4be7-b292-29b7b9adc7fb.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d5674dff41a957be50bd534a7ba6fc1c5940af21/e2e/65a8a149-dbf1-4be7-b292-29b7b9adc7fb.md."

$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d5674dff41a957be50bd534a7ba6fc1c5940af21/e2e/65a8a149-dbf1-4be7-b292-29b7b9adc7fb.md"
$mdDisplay = "65a8a149-dbf1-4be7-b292-29b7b9adc7fb.md"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

# Widen the "Error Detail" column (P / col 16) so the long message fits.
$ws2.Columns.Item(16).ColumnWidth = 39.17

# Row 8 is the 65a8a149-... entry.
# I8 = Latest Target File -> hyperlink to the handback markdown file.
$ws2.Hyperlinks.Add($ws2.Cells.Item(8, 9), $mdUrl, "", "", $mdDisplay)

# J8 = Latest Handback File
$ws2.Range("J8").Value = "65a8a149-dbf1-4be7-b292-29b7b9adc7fb.068649dc21e79ad254d63a4b3512f23e2c2171f4.zh-cn.xlf"

# K8 = Latest Handback DateTime
$ws2.Range("K8").Value = "2016-09-04 22:47:30"

# P8 = Error Detail
$ws2.Range("P8").Value = $errorDetail

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

# Widen the "Error Detail" column (P / col 16) so the long message fits.
$ws3.Columns.Item(16).ColumnWidth = 39.17

# Row 8 is the 65a8a149-... entry.
# I8 = Latest Target File -> hyperlink to the handback markdown file.
$ws3.Hyperlinks.Add($ws3.Cells.Item(8, 9), $mdUrl, "", "", $mdDisplay)

# J8 = Latest Handback File
$ws3.Range("J8").Value = "65a8a149-dbf1-4be7-b292-29b7b9adc7fb.068649dc21e79ad254d63a4b3512f23e2c2171f4.de-de.xlf"

# K8 = Latest Handback DateTime
$ws3.Range("K8").Value = "2016-09-04 22:47:37"

# P8 = Error Detail
$ws3.Range("P8").Value = $errorDetail
